# chore: fix all test rules
# Update Risk_Label (column H) values from "No Match" to the correct
# computed risk levels for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16  = "Medium"
    23  = "Low"
    37  = "Low"
    47  = "Low"
    93  = "Low"
    102 = "Low"
    104 = "Low"
    105 = "Low"
    110 = "Low"
    111 = "Low"
    113 = "Low"
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}
